$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new item row for "CONTAFEVER N 200MG/5ML SUSP. 120ML" right
#    after row 7 (ALKAPRESS) and before the old row 8 (DECLOPHEN), keeping
#    the list sorted alphabetically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "'CONTAFEVER N 200MG/5ML SUSP. 120ML"
$ws.Range("H8").Value = "'20:0"
$ws.Range("L8").Value = "'1"
$ws.Range("N8").Value = "'33.00"
$ws.Range("P8").Value = "'33.0000"
$ws.Range("Q8").Value = "'1:0"

# Re-apply the row-7 formatting (styles + merges) onto the freshly inserted
# row 8 without disturbing the text values just written above.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Renumber the items that used to be #2-#8 (now #3-#9) in column A.
#    (They are plain numbers, not formulas, so Excel does not auto-renumber
#    them on insert.)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8

# ---------------------------------------------------------------------------
# 3) Insert a new item row for "VOLTAREN 75MG/3ML 3 AMP." right after the
#    old OMEZ row (now row 14) and before the سرنجات row (now row 15).
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "'VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H15").Value = "'0:2"
$ws.Range("L15").Value = "'1"
$ws.Range("N15").Value = "'51.00"
$ws.Range("P15").Value = "'33.6600"
$ws.Range("Q15").Value = "'0:2"

$ws.Range("A14:Q14").Copy()
$ws.Range("A15:Q15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) سرنجات 3 سم is now row 16: bump its item number and refresh the
#    "سعر البيع" / "عدد التعاملات" figures.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 10
$ws.Range("P16").Value = "'18.0000"
$ws.Range("Q16").Value = "'9:0"

$ws.Range("A16:Q16").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) Update the grand-total cell (old row 15, now row 17) to match the new
#    sum of the "سعر البيع" column across the (now) 10 item rows.
# ---------------------------------------------------------------------------
$ws.Range("P17").Value = 332.76

# ---------------------------------------------------------------------------
# 6) Refresh the generation timestamp in the footer (old row 16, now row 18).
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "Saturday, 30 August, 2025 11:11 AM"

$excel.CutCopyMode = 0
